# Update the date line
$d = $word.ActiveDocument

$pairs = @(
    @("2024-02-26 Monday", "2024-02-27 Tuesday"),
    @("261÷7=", "812÷3="),
    @("469÷3=", "761÷4="),
    @("794÷6=", "533÷4="),
    @("514÷2=", "179÷8="),
    @("771÷6=", "244÷2="),
    @("514÷7=", "797÷5="),
    @("971÷7=", "659÷2="),
    @("196÷3=", "904÷4="),
    @("357÷8=", "361÷7="),
    @("591÷2=", "646÷8="),
    @("731÷3=", "414÷9="),
    @("283÷6=", "935÷9="),
    @("404÷7=", "418÷4="),
    @("106÷6=", "314÷8="),
    @("180÷5=", "284÷6="),
    @("104÷4=", "384÷7="),
    @("503÷6=", "438÷2="),
    @("368÷8=", "531÷8="),
    @("275÷9=", "367÷3="),
    @("476÷9=", "982÷7="),
    @("536÷4=", "738÷6="),
    @("720÷9=", "502÷8="),
    @("840÷8=", "860÷3="),
    @("792÷6=", "417÷8="),
    @("578÷9=", "422÷7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
